# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the "Estado de Cuenta" detail table (rows 16-33) so that instead
# of being grouped first by worker (all periods for MARTIN, then all
# periods for PATRICIA), the rows alternate worker-by-worker within each
# period, listed in ascending period order (2009 .. 2105). Also updates the
# "Valor Mora" amounts for the newly-reordered rows to match the refreshed
# account-statement data (periods 2009-2104 now carry 35112, and the final
# 2105 pair carries 28090).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rowsData = @{
    16 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2009"; F = 35112 }
    17 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2009"; F = 35112 }
    18 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2010"; F = 35112 }
    19 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2010"; F = 35112 }
    20 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2011"; F = 35112 }
    21 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2011"; F = 35112 }
    22 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2012"; F = 35112 }
    23 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2012"; F = 35112 }
    24 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2101"; F = 35112 }
    25 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2101"; F = 35112 }
    26 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2102"; F = 35112 }
    27 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2102"; F = 35112 }
    28 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2103"; F = 35112 }
    29 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2103"; F = 35112 }
    30 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2104"; F = 35112 }
    31 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2104"; F = 35112 }
    32 = @{ C = "73113375"; D = "MARTIN EMILIO ROSALES GOMEZ";      E = "2105"; F = 28090 }
    33 = @{ C = "45468832"; D = "PATRICIA DEL CARMEN ROSALES GOMEZ"; E = "2105"; F = 28090 }
}

foreach ($r in $rowsData.Keys) {
    $row = $rowsData[$r]
    # Column B (Tipo Doc Trabajador = "CC") is unchanged for every row.
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
}
